# Add the MCH200 collection metadata as row 2, directly under the existing
# header row (identifier | alternativeIdentifiers | title | date_s |
# levelOfDescription | extentAndMedium | notes | file_path).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# identifier
$ws.Range("A2").Value = "MCH200"
# alternativeIdentifiers (B2) is intentionally left untouched / blank.
# title
$ws.Range("C2").Value = "NEWSLETTERS OF THE ANTI APARTHEID MOVEMENT 1989-1990"
# date_s
$ws.Range("D2").Value = "1989-1990"
# levelOfDescription
$ws.Range("E2").Value = "Series"
# extentAndMedium
$ws.Range("F2").Value = "1 Box"
# notes
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"
# file_path (H2) carries no value, only the row's formatting.

# Give the new row the same look as the rest of the data rows: 10pt Calibri,
# automatic/theme text color (as opposed to the bold header style above it).
$dataFont = $ws.Range("A2").Font
$dataFont.Name = "Calibri"
$dataFont.ThemeColor = 1

# Copy that formatting across the rest of the populated row (including the
# empty file_path cell) so every cell in row 2 shares one consistent style.
$ws.Range("A2").Copy()
$ws.Range("C2:H2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Re-select the newly entered row and keep the header row frozen, matching
# the workbook's saved view state.
[void]$ws.Range("A2:J2").Select()
$excel.ActiveWindow.FreezePanes = $true
